# working on latitudinal gradient
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: add a new "new" column (F) with data, and re-style the last data
# row (Flyvestation, DA / row 10) with a red font.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("F1").Value = "new"

$sheet1F = @(5, 4, 10, 4, 3, 5, 2, 0, 0, 0)
for ($i = 0; $i -lt $sheet1F.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $sheet1F[$i]
}

# Row 10 (Flyvestation, DA) gets a red font applied across A:F.
$ws1.Range("A10:F10").Font.Color = 255

$ws1.Range("E20").Select()

# ---------------------------------------------------------------------------
# Sheet2: view changes only (scroll position + selection).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A17").Select()
$ws2.Application.ActiveWindow.ScrollRow = 17
$ws2.Range("E24").Select()

# ---------------------------------------------------------------------------
# Sheet3: add a new "new" column (F) with data.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("F1").Value = "new"

$sheet3F = @(6, 14, 6, 10, 8, 17, 7, 11)
for ($i = 0; $i -lt $sheet3F.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 6).Value = $sheet3F[$i]
}

$ws3.Range("F9").Select()
$ws3.Activate()

# ---------------------------------------------------------------------------
# Sheet4: no longer the active/selected tab (Sheet3 is now active instead).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A34").Select()

# Final active sheet/tab is Sheet3 (activeTab index 2, 0-based).
$ws3.Activate()
